# globaladmin added as default in zone-user tables.
# Adds two default zone-user rows (eng/MOR and fra/MOR) for user "globaladmin".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$boolFormat = """TRUE"";""TRUE"";""FALSE"""

# Row 2: eng / MOR / globaladmin / TRUE / now()
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "MOR"
$ws.Range("C2").Value = "globaladmin"
$ws.Range("D2").Value = $true
$ws.Range("D2").NumberFormat = $boolFormat
$ws.Range("E2").Value = "now()"

# Row 3: fra / MOR / globaladmin / TRUE / now()
$ws.Range("A3").Value = "fra"
$ws.Range("B3").Value = "MOR"
$ws.Range("C3").Value = "globaladmin"
$ws.Range("D3").Value = $true
$ws.Range("D3").NumberFormat = $boolFormat
$ws.Range("E3").Value = "now()"

# Match the authored column widths for zone_code / usr_id / eff_dtimes columns.
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 14.1
$ws.Columns.Item(5).ColumnWidth = 13.8

# Move the active selection past the newly-entered rows.
[void]$ws.Range("A4").Select()
